# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates to the Kujata_Profits workbook
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 5128.643
$ws.Range("I76").Value = 5225.375
$ws.Range("K76").Value = 5225.375
$ws.Range("M76").Value = -4910.375
# Row 79
$ws.Range("H79").Value = 5128.643
$ws.Range("I79").Value = 5225.375
$ws.Range("K79").Value = 5225.375
$ws.Range("M79").Value = -4133.375
# Row 92
$ws.Range("H92").Value = 788.9231
$ws.Range("I92").Value = 419.5
$ws.Range("J92").Value = 1620.125
$ws.Range("K92").Value = 419.5
$ws.Range("L92").Value = 1620.125
$ws.Range("M92").Value = 828.5
$ws.Range("N92").Value = -4116.125
# Row 121
$ws.Range("H121").Value = 527.8570999999999
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 527.8570999999999
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 1583.5713
$ws.Range("M121").Value = $null
$ws.Range("N121").Value = -5077.5713
# Row 137
$ws.Range("H137").Value = 1416.0834
$ws.Range("I137").Value = 1118.875
$ws.Range("J137").Value = 2010.5
$ws.Range("K137").Value = 3356.625
$ws.Range("L137").Value = 6031.5
$ws.Range("M137").Value = -806.625
$ws.Range("N137").Value = -11131.5
# Row 138
$ws.Range("H138").Value = 445479.38
$ws.Range("J138").Value = 545322.7
$ws.Range("L138").Value = 1635968.1
$ws.Range("N138").Value = -1646248.1
# Row 141
$ws.Range("H141").Value = 3188.2
$ws.Range("I141").Value = 3188.2
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9564.599999999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -4384.599999999999
$ws.Range("N141").Value = $null

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6432.14
$ws.Range("I32").Value = 5267.244
$ws.Range("K32").Value = 5267.244
$ws.Range("M32").Value = -4980.244
# Row 61
$ws.Range("H61").Value = 62501290
$ws.Range("I61").Value = 76924160
$ws.Range("J61").Value = 2171.3333
$ws.Range("K61").Value = 76924160
$ws.Range("L61").Value = 2171.3333
$ws.Range("M61").Value = -76923948
$ws.Range("N61").Value = -2595.3333
# Row 74
$ws.Range("H74").Value = 2481.2222
$ws.Range("I74").Value = 1762
$ws.Range("J74").Value = 2938.9092
$ws.Range("K74").Value = 1762
$ws.Range("L74").Value = 2938.9092
$ws.Range("M74").Value = -888
$ws.Range("N74").Value = -4686.9092
# Row 77
$ws.Range("H77").Value = 2481.2222
$ws.Range("I77").Value = 1762
$ws.Range("J77").Value = 2938.9092
$ws.Range("K77").Value = 8810
$ws.Range("L77").Value = 14694.546
$ws.Range("M77").Value = -4442
$ws.Range("N77").Value = -23430.546
# Row 122
$ws.Range("H122").Value = 1933.75
$ws.Range("I122").Value = 1546.75
$ws.Range("J122").Value = 3481.75
$ws.Range("K122").Value = 4640.25
$ws.Range("L122").Value = 10445.25
$ws.Range("M122").Value = -2190.25
$ws.Range("N122").Value = -15345.25
# Row 132
$ws.Range("H132").Value = 2432.6135
$ws.Range("I132").Value = 2028.9259
$ws.Range("J132").Value = 3073.7646
$ws.Range("K132").Value = 6086.7777
$ws.Range("L132").Value = 9221.293799999999
$ws.Range("M132").Value = -3556.7777
$ws.Range("N132").Value = -14281.2938
# Row 136
$ws.Range("H136").Value = 62501290
$ws.Range("I136").Value = 76924160
$ws.Range("J136").Value = 2171.3333
$ws.Range("K136").Value = 230772480
$ws.Range("L136").Value = 6513.999899999999
$ws.Range("M136").Value = -230769930
$ws.Range("N136").Value = -11613.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3299.9167
$ws.Range("I20").Value = 3362.375
$ws.Range("K20").Value = 3362.375
$ws.Range("M20").Value = -3115.375
# Row 86
$ws.Range("H86").Value = 2516.0571
$ws.Range("I86").Value = 2422.4546
$ws.Range("J86").Value = 2674.4614
$ws.Range("K86").Value = 2422.4546
$ws.Range("L86").Value = 2674.4614
$ws.Range("M86").Value = -1299.4546
$ws.Range("N86").Value = -4920.4614
# Row 89
$ws.Range("H89").Value = 2516.0571
$ws.Range("I89").Value = 2422.4546
$ws.Range("J89").Value = 2674.4614
$ws.Range("K89").Value = 12112.273
$ws.Range("L89").Value = 13372.307
$ws.Range("M89").Value = -6496.273000000001
$ws.Range("N89").Value = -24604.307
# Row 105
$ws.Range("H105").Value = 250001740
$ws.Range("I105").Value = 500000500
$ws.Range("K105").Value = 500000500
$ws.Range("M105").Value = -499998753
# Row 107
$ws.Range("H107").Value = 1220.1765
$ws.Range("I107").Value = 920.9091
$ws.Range("J107").Value = 1768.8334
$ws.Range("K107").Value = 920.9091
$ws.Range("L107").Value = 1768.8334
$ws.Range("M107").Value = 999.0909
$ws.Range("N107").Value = -5608.8334
# Row 134
$ws.Range("H134").Value = 1660.5
$ws.Range("I134").Value = 1469.1428
$ws.Range("K134").Value = 4407.428400000001
$ws.Range("M134").Value = -1872.428400000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1444.683
$ws.Range("I31").Value = 1350.8
$ws.Range("J31").Value = 5200
$ws.Range("K31").Value = 1350.8
$ws.Range("L31").Value = 5200
$ws.Range("M31").Value = -1055.8
$ws.Range("N31").Value = -5790
# Row 34
$ws.Range("H34").Value = 1444.683
$ws.Range("I34").Value = 1350.8
$ws.Range("J34").Value = 5200
$ws.Range("K34").Value = 1350.8
$ws.Range("L34").Value = 5200
$ws.Range("M34").Value = -1148.8
$ws.Range("N34").Value = -5604
# Row 132
$ws.Range("H132").Value = 2572.077
$ws.Range("I132").Value = 1867.375
$ws.Range("J132").Value = 3699.6
$ws.Range("K132").Value = 5602.125
$ws.Range("L132").Value = 11098.8
$ws.Range("M132").Value = -3072.125
$ws.Range("N132").Value = -16158.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 692.4783
$ws.Range("J5").Value = 558.75
$ws.Range("L5").Value = 1676.25
$ws.Range("N5").Value = -1900.25
# Row 12
$ws.Range("H12").Value = 116.5
$ws.Range("J12").Value = 78.066666
$ws.Range("L12").Value = 234.199998
$ws.Range("N12").Value = -580.1999980000001
# Row 26
$ws.Range("H26").Value = 91.40000000000001
$ws.Range("I26").Value = 66.28570999999999
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 198.85713
$ws.Range("L26").Value = 450
$ws.Range("M26").Value = 89.14287000000002
$ws.Range("N26").Value = -1026
# Row 131
$ws.Range("H131").Value = 27030470
$ws.Range("J131").Value = 4875.44
$ws.Range("L131").Value = 14626.32
$ws.Range("N131").Value = -24706.32
# Row 135
$ws.Range("H135").Value = 692.4783
$ws.Range("J135").Value = 558.75
$ws.Range("L135").Value = 5028.75
$ws.Range("N135").Value = -10098.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 34618750
$ws.Range("I70").Value = 27781362
$ws.Range("J70").Value = 50002876
$ws.Range("K70").Value = 27781362
$ws.Range("L70").Value = 50002876
$ws.Range("M70").Value = -27781092
$ws.Range("N70").Value = -50003416
# Row 73
$ws.Range("H73").Value = 34618750
$ws.Range("I73").Value = 27781362
$ws.Range("J73").Value = 50002876
$ws.Range("K73").Value = 27781362
$ws.Range("L73").Value = 50002876
$ws.Range("M73").Value = -27780426
$ws.Range("N73").Value = -50004748
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null
# Row 126
$ws.Range("H126").Value = 1868.45
$ws.Range("I126").Value = 1646.2
$ws.Range("J126").Value = 2090.7
$ws.Range("K126").Value = 4938.6
$ws.Range("L126").Value = 6272.099999999999
$ws.Range("M126").Value = -2468.6
$ws.Range("N126").Value = -11212.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3926.318
$ws.Range("I40").Value = 1692.0834
$ws.Range("K40").Value = 1692.0834
$ws.Range("M40").Value = -1556.0834
# Row 68
$ws.Range("H68").Value = 1287.2858
$ws.Range("I68").Value = 1281.2
$ws.Range("K68").Value = 1281.2
$ws.Range("M68").Value = -532.2
# Row 71
$ws.Range("H71").Value = 1287.2858
$ws.Range("I71").Value = 1281.2
$ws.Range("K71").Value = 6406
$ws.Range("M71").Value = -2662

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Range("H104").Value = 23380
$ws.Range("J104").Value = 23380
$ws.Range("L104").Value = 23380
$ws.Range("N104").Value = -30368
# Row 107
$ws.Range("H107").Value = 397.8889
$ws.Range("I107").Value = 298.06668
$ws.Range("J107").Value = 522.6667
$ws.Range("K107").Value = 894.2000400000001
$ws.Range("L107").Value = 1568.0001
$ws.Range("M107").Value = 1025.79996
$ws.Range("N107").Value = -5408.0001
# Row 126
$ws.Range("H126").Value = 43479132
$ws.Range("I126").Value = 62500476
$ws.Range("K126").Value = 187501428
$ws.Range("M126").Value = -187498958
